$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 29.32133366666666
$ws.Range("H2").Value = 87.964001
$ws.Range("I2").Value = 0.006401919837078288
$ws.Range("J2").Value = 0.006401919837078288
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 17.88428466666667
$ws.Range("N2").Value = 53.652854
$ws.Range("O2").Value = 0.4435785307770658
$ws.Range("P2").Value = 0.4435785307770658
$ws.Range("Q2").Value = 524.3910781009838
$ws.Range("R2").Value = 4719.519702908854
$ws.Range("S2").Value = 0.002839754195483739
$ws.Range("T2").Value = 0.00283975419548374

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 29.32133366666666
$ws.Range("H3").Value = 87.964001
$ws.Range("I3").Value = 0.006401919837078288
$ws.Range("J3").Value = 0.006401919837078288
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 12.393653
$ws.Range("N3").Value = 37.180959
$ws.Range("O3").Value = 0.3073960458115111
$ws.Range("P3").Value = 0.3073960458115112
$ws.Range("Q3").Value = 363.3984349618843
$ws.Range("R3").Value = 3270.585914656959
$ws.Range("S3").Value = 0.001967924843520139
$ws.Range("T3").Value = 0.00196792484352014

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 29.32133366666666
$ws.Range("H4").Value = 87.964001
$ws.Range("I4").Value = 0.006401919837078288
$ws.Range("J4").Value = 0.006401919837078288
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 5.643036666666666
$ws.Range("N4").Value = 16.92911
$ws.Range("O4").Value = 0.1399625403182342
$ws.Range("P4").Value = 0.1399625403182342
$ws.Range("Q4").Value = 165.4613609965677
$ws.Range("R4").Value = 1489.15224896911
$ws.Range("S4").Value = 0.0008960289633111729
$ws.Range("T4").Value = 0.000896028963311173

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 29.32133366666666
$ws.Range("H5").Value = 87.964001
$ws.Range("I5").Value = 0.006401919837078288
$ws.Range("J5").Value = 0.006401919837078288
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 4.397218333333333
$ws.Range("N5").Value = 13.191655
$ws.Range("O5").Value = 0.1090628830931889
$ws.Range("P5").Value = 0.1090628830931889
$ws.Range("Q5").Value = 128.9323059568505
$ws.Range("R5").Value = 1160.390753611655
$ws.Range("S5").Value = 0.0006982118347632362
$ws.Range("T5").Value = 0.0006982118347632364

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 45.524413
$ws.Range("H6").Value = 136.573239
$ws.Range("I6").Value = 0.009939644832300594
$ws.Range("J6").Value = 0.009939644832300592
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 17.88428466666667
$ws.Range("N6").Value = 53.652854
$ws.Range("O6").Value = 0.4435785307770658
$ws.Range("P6").Value = 0.4435785307770658
$ws.Range("Q6").Value = 814.1715613749009
$ws.Range("R6").Value = 7327.544052374106
$ws.Range("S6").Value = 0.004409013051157752
$ws.Range("T6").Value = 0.004409013051157751

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 45.524413
$ws.Range("H7").Value = 136.573239
$ws.Range("I7").Value = 0.009939644832300594
$ws.Range("J7").Value = 0.009939644832300592
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 12.393653
$ws.Range("N7").Value = 37.180959
$ws.Range("O7").Value = 0.3073960458115111
$ws.Range("P7").Value = 0.3073960458115112
$ws.Range("Q7").Value = 564.213777750689
$ws.Range("R7").Value = 5077.923999756201
$ws.Range("S7").Value = 0.003055407518220023
$ws.Range("T7").Value = 0.003055407518220023

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 45.524413
$ws.Range("H8").Value = 136.573239
$ws.Range("I8").Value = 0.009939644832300594
$ws.Range("J8").Value = 0.009939644832300592
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 5.643036666666666
$ws.Range("N8").Value = 16.92911
$ws.Range("O8").Value = 0.1399625403182342
$ws.Range("P8").Value = 0.1399625403182342
$ws.Range("Q8").Value = 256.8959317874767
$ws.Range("R8").Value = 2312.06338608729
$ws.Range("S8").Value = 0.0013911779405898
$ws.Range("T8").Value = 0.001391177940589799

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 45.524413
$ws.Range("H9").Value = 136.573239
$ws.Range("I9").Value = 0.009939644832300594
$ws.Range("J9").Value = 0.009939644832300592
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 4.397218333333333
$ws.Range("N9").Value = 13.191655
$ws.Range("O9").Value = 0.1090628830931889
$ws.Range("P9").Value = 0.1090628830931889
$ws.Range("Q9").Value = 200.1807834578383
$ws.Range("R9").Value = 1801.627051120545
$ws.Range("S9").Value = 0.001084046322333019
$ws.Range("T9").Value = 0.001084046322333019

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 4438.215250666667
$ws.Range("H10").Value = 13314.645752
$ws.Range("I10").Value = 0.9690247577915309
$ws.Range("J10").Value = 0.9690247577915307
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 17.88428466666667
$ws.Range("N10").Value = 53.652854
$ws.Range("O10").Value = 0.4435785307770658
$ws.Range("P10").Value = 0.4435785307770658
$ws.Range("Q10").Value = 79374.30495486404
$ws.Range("R10").Value = 714368.7445937763
$ws.Range("S10").Value = 0.4298385783477693
$ws.Range("T10").Value = 0.4298385783477693

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 4438.215250666667
$ws.Range("H11").Value = 13314.645752
$ws.Range("I11").Value = 0.9690247577915309
$ws.Range("J11").Value = 0.9690247577915307
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 12.393653
$ws.Range("N11").Value = 37.180959
$ws.Range("O11").Value = 0.3073960458115111
$ws.Range("P11").Value = 0.3073960458115112
$ws.Range("Q11").Value = 55005.6997560707
$ws.Range("R11").Value = 495051.2978046362
$ws.Range("S11").Value = 0.2978743788385739
$ws.Range("T11").Value = 0.2978743788385739

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 4438.215250666667
$ws.Range("H12").Value = 13314.645752
$ws.Range("I12").Value = 0.9690247577915309
$ws.Range("J12").Value = 0.9690247577915307
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 5.643036666666666
$ws.Range("N12").Value = 16.92911
$ws.Range("O12").Value = 0.1399625403182342
$ws.Range("P12").Value = 0.1399625403182342
$ws.Range("Q12").Value = 25045.01139407119
$ws.Range("R12").Value = 225405.1025466407
$ws.Range("S12").Value = 0.1356271667317642
$ws.Range("T12").Value = 0.1356271667317642

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 4438.215250666667
$ws.Range("H13").Value = 13314.645752
$ws.Range("I13").Value = 0.9690247577915309
$ws.Range("J13").Value = 0.9690247577915307
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 4.397218333333333
$ws.Range("N13").Value = 13.191655
$ws.Range("O13").Value = 0.1090628830931889
$ws.Range("P13").Value = 0.1090628830931889
$ws.Range("Q13").Value = 19515.80146751106
$ws.Range("R13").Value = 175642.2132075996
$ws.Range("S13").Value = 0.1056846338734234
$ws.Range("T13").Value = 0.1056846338734234

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 67.02347933333333
$ws.Range("H14").Value = 201.070438
$ws.Range("I14").Value = 0.01463367753909034
$ws.Range("J14").Value = 0.01463367753909034
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 17.88428466666667
$ws.Range("N14").Value = 53.652854
$ws.Range("O14").Value = 0.4435785307770658
$ws.Range("P14").Value = 0.4435785307770658
$ws.Range("Q14").Value = 1198.666983747784
$ws.Range("R14").Value = 10788.00285373005
$ws.Range("S14").Value = 0.00649118518265504
$ws.Range("T14").Value = 0.006491185182655041

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 67.02347933333333
$ws.Range("H15").Value = 201.070438
$ws.Range("I15").Value = 0.01463367753909034
$ws.Range("J15").Value = 0.01463367753909034
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 12.393653
$ws.Range("N15").Value = 37.180959
$ws.Range("O15").Value = 0.3073960458115111
$ws.Range("P15").Value = 0.3073960458115112
$ws.Range("Q15").Value = 830.6657457100047
$ws.Range("R15").Value = 7475.991711390042
$ws.Range("S15").Value = 0.004498334611197095
$ws.Range("T15").Value = 0.004498334611197096

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 67.02347933333333
$ws.Range("H16").Value = 201.070438
$ws.Range("I16").Value = 0.01463367753909034
$ws.Range("J16").Value = 0.01463367753909034
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 5.643036666666666
$ws.Range("N16").Value = 16.92911
$ws.Range("O16").Value = 0.1399625403182342
$ws.Range("P16").Value = 0.1399625403182342
$ws.Range("Q16").Value = 378.2159514055755
$ws.Range("R16").Value = 3403.94356265018
$ws.Range("S16").Value = 0.002048166682568969
$ws.Range("T16").Value = 0.002048166682568969

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 67.02347933333333
$ws.Range("H17").Value = 201.070438
$ws.Range("I17").Value = 0.01463367753909034
$ws.Range("J17").Value = 0.01463367753909034
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 4.397218333333333
$ws.Range("N17").Value = 13.191655
$ws.Range("O17").Value = 0.1090628830931889
$ws.Range("P17").Value = 0.1090628830931889
$ws.Range("Q17").Value = 294.7168720883211
$ws.Range("R17").Value = 2652.45184879489
$ws.Range("S17").Value = 0.001595991062669234
$ws.Range("T17").Value = 0.001595991062669234
